$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three rows that only held the instructors' names in columns B/C
# (with no label in column A) are removed; everything below shifts up.
$ws.Rows("13:15").Delete()

# After the shift, several label rows need their B/C content updated so the
# text lines up with the (now one-row-higher) label it sits beside.
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B15").Value = "5840793 - Sérgio Schneider"
$ws.Range("C15").Value = "5840793 - Sérgio Schneider"

$ws.Range("B18").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C18").Value = "7797767 - Viktor Pastoukhov"

$ws.Range("B19").Value = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"
$ws.Range("C19").Value = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"

$ws.Range("B20").Value = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("C20").Value = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."

$ws.Range("B21").Value = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("C21").Value = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
